# Daily attendance processing - 2025-10-23 12:39:31
# Applies the session-analysis refresh: recorder-list reordering, updated
# attendance counts/percentages, and the "Year 2 / A2 / HISTOLOGY" session
# (row 22) flipping from Pending/Not-Recorded (pink) to Recorded (green).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Year 2 / A1 / ANATOMY 23/10) -----------------------------------
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("H2").Value = "95/216"

# --- Class statistics: Recorded Sessions / Missing Sessions ----------------
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 1

# --- Row 7 (Year 2 / A1 / HISTOLOGY 23/10) ----------------------------------
$ws.Range("G7").Value = "Arwa.elnagar@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("H7").Value = "189/216"

# --- Class statistics: Coverage % / Average Attendance % -------------------
# (these columns store percentages as literal text, e.g. "10.8%" -- force
# Text number format first so the COM layer doesn't coerce the string into
# a numeric percentage value)
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "10.8%"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "46.4%"

# --- Group statistics row 15/16 (Year 2 / A1, A2) ---------------------------
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "65.7%"
$ws.Range("O16").Value = 2
$ws.Range("P16").Value = 0
$ws.Range("R16").NumberFormat = "@"
$ws.Range("R16").Value = "13.3%"
$ws.Range("S16").NumberFormat = "@"
$ws.Range("S16").Value = "22.1%"

# --- Row 17 (Year 2 / A2 / ANATOMY 23/10) -----------------------------------
$ws.Range("G17").Value = "servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("H17").Value = "95/217"

# --- Row 22 (Year 2 / A2 / HISTOLOGY 23/10): now recorded, fill -> green ----
$ws.Range("A22:I22").Interior.Color = 9498256
$ws.Range("G22").Value = "Arwa.elnagar@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("H22").Value = "1/217"
$ws.Range("I22").Value = "Recorded"

# --- Recorder-list reordering (same people, canonical order updated) -------
$ws.Range("G32").Value = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G47").Value = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G62").Value = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G77").Value = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G92").Value = "servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G107").Value = "servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G112").Value = "Arwa.elnagar@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
